$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.885.55'
$ws.Range('E2').Value = '  +0.83%  '

$ws.Range('D3').Value = '1.709.56'
$ws.Range('E3').Value = '  +0.86%  '

$ws.Range('D4').Value = "'1.013"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +1.28%  '

$ws.Range('D5').Value = "'318.31"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.33%  '

$ws.Range('D6').Value = "'1.012"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.16%  '

$ws.Range('D7').Value = "'0.3964"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').Value = "'0.4109"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.82%  '

$ws.Range('D9').Value = "'1.524"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.94%  '

$ws.Range('D10').Value = "'1.013"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.29%  '

$ws.Range('D11').Value = "'52.28"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.04%  '

$ws.Range('D12').Value = "'0.08853"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.72%  '

$ws.Range('D13').Value = "'7.671"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.49%  '

$ws.Range('D14').Value = "'24.66"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.41%  '

$ws.Range('D15').Value = "'0.00001383"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.23%  '

$ws.Range('D16').Value = "'8.098"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.09%  '

$ws.Range('D17').Value = '1.713.00'
$ws.Range('E17').Value = '  +1.20%  '

$ws.Range('D18').Value = "'100.41"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.08%  '

$ws.Range('D19').Value = "'0.07129"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.68%  '

$ws.Range('D20').Value = "'20.07"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.65%  '

$ws.Range('D21').Value = "'7.446"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.73%  '

$ws.Range('D22').Value = "'1.013"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.28%  '

$ws.Range('D23').Value = "'14.42"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.77%  '

$ws.Range('D24').Value = '24.883.86'
$ws.Range('E24').Value = '  +0.81%  '

$ws.Range('D25').Value = "'3.053"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.60%  '

$ws.Range('D26').Value = "'2.358"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.59%  '

$ws.Range('D27').Value = "'22.95"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.07%  '

$ws.Range('D28').Value = "'164.95"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.98%  '

$ws.Range('D29').Value = "'8.685"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +15.05%  '

$ws.Range('D30').Value = "'139.18"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.56%  '

$ws.Range('D31').Value = "'5.205"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.45%  '

$ws.Range('B32').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C32').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D32').Value = '1.901.66'
$ws.Range('E32').Value = '  +1.06%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'0.09010"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.50%  '

$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = "'7.656"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.81%  '

$ws.Range('D35').Value = "'1.059"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.48%  '

$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = "'1.992"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.69%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.02954"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.12%  '

$ws.Range('D38').Value = "'0.2759"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.43%  '

$ws.Range('D39').Value = "'10.91"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.19%  '

$ws.Range('D40').Value = "'14.52"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.22%  '

$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = "'0.09220"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.32%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = "'0.7974"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.62%  '

$ws.Range('D43').Value = "'1.482"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.28%  '

$ws.Range('D44').Value = "'16.58"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.37%  '

$ws.Range('D45').Value = "'0.7315"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.54%  '

$ws.Range('D46').Value = "'2.626"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.26%  '

$ws.Range('D47').Value = "'4.275"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.00%  '

$ws.Range('D48').Value = "'1.011"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.01%  '

$ws.Range('D49').Value = "'1.335"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.51%  '

$ws.Range('E50').Value = '  -0.38%  '

$ws.Range('D51').Value = "'92.24"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.36%  '
